# Enable non-blocking parallel processing of ticket codes
#
# The "Daten" sheet has a stray leftover cell (T25 = "sni") far outside the
# real data table (A1:F13). Clearing it shrinks the sheet's used range back
# down to A1:T13, matching the cleaned-up workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daten")

$ws.Range("T25").ClearContents()
